# Add a new "Thank you" closing slide at the end of the deck.
$p = $ppt.ActivePresentation

# Append a new slide (index 3) using the "Title and Content" layout,
# the same layout slide 1 / other content slides use (slideLayout2.xml,
# ppLayoutText = 2).
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

# Title placeholder -> "Thank you"
$titleShape = $newSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Thank you"
$titleShape.Name = "제목 1"

# Content placeholder stays empty, just rename to match authoring locale.
$bodyShape = $newSlide.Shapes.Item(2)
$bodyShape.Name = "내용 개체 틀 2"
